# Weekly fruit/vegetable data update.
# A new weekly price record for "Ají" (Inferno, Primera) at Feria Lagunitas
# de Puerto Montt needs to be inserted right after the existing row for
# date 44777 (row 467), which pushes all subsequent rows down by one and
# extends the used range from A1:R495 to A1:R496.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 467 - everything currently at/after row 467
# shifts down one row (467 -> 468, ..., 495 -> 496).
$ws.Rows("467:467").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(467, 1).Value = 4
$ws.Cells.Item(467, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(467, 3).Value = "Los Lagos"
$ws.Cells.Item(467, 4).Value = 45265
$ws.Cells.Item(467, 5).Value = 10
$ws.Cells.Item(467, 6).Value = 100112021
$ws.Cells.Item(467, 7).Value = "Ají"
$ws.Cells.Item(467, 8).Value = "Inferno"
$ws.Cells.Item(467, 9).Value = "Primera"
$ws.Cells.Item(467, 10).Value = 180
$ws.Cells.Item(467, 11).Value = 45000
$ws.Cells.Item(467, 12).Value = 45000
$ws.Cells.Item(467, 13).Value = 45000
$ws.Cells.Item(467, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(467, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(467, 16).Value = 4500
$ws.Cells.Item(467, 17).Value = 10
$ws.Cells.Item(467, 18).Value = "Hortaliza"
